$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsZhCn.Range("D5").Value = "2016-03-04 02:59:02"
$wsDeDe.Range("D5").Value = "2016-03-04 02:59:15"
